# edit.ps1
# Applies the CasosColombia.xlsx update:
#  - Retroactive corrections to column BQ (rows 18-187)
#  - A handful of scattered numeric<->"NaN" corrections (I22, CU26, DT34, AP.., CF93/CF108)
#  - Append a new data row (188) for 2020-09-08 with a full set of columns A:DX
#  - Move the active selection to the new last cell (DX188)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Column BQ (rows 18-187): updated running totals
# ---------------------------------------------------------------------------
$bqUpdates = @(
    @(18, 1),
    @(19, 1),
    @(20, 2),
    @(21, 2),
    @(22, 2),
    @(23, 2),
    @(24, 3),
    @(25, 4),
    @(26, 4),
    @(27, 11),
    @(28, 16),
    @(29, 16),
    @(30, 16),
    @(31, 16),
    @(32, 16),
    @(33, 16),
    @(34, 16),
    @(35, 17),
    @(36, 18),
    @(37, 19),
    @(38, 31),
    @(39, 31),
    @(40, 31),
    @(41, 31),
    @(42, 31),
    @(43, 31),
    @(44, 31),
    @(45, 31),
    @(46, 31),
    @(47, 33),
    @(48, 33),
    @(49, 33),
    @(50, 33),
    @(51, 37),
    @(52, 43),
    @(53, 43),
    @(54, 43),
    @(55, 43),
    @(56, 43),
    @(57, 43),
    @(58, 43),
    @(59, 45),
    @(60, 48),
    @(61, 48),
    @(62, 49),
    @(63, 49),
    @(64, 49),
    @(65, 49),
    @(66, 50),
    @(67, 51),
    @(68, 51),
    @(69, 51),
    @(70, 51),
    @(71, 51),
    @(72, 51),
    @(73, 54),
    @(74, 55),
    @(75, 56),
    @(76, 56),
    @(77, 56),
    @(78, 57),
    @(79, 57),
    @(80, 61),
    @(81, 63),
    @(82, 63),
    @(83, 64),
    @(84, 77),
    @(85, 90),
    @(86, 110),
    @(87, 154),
    @(88, 188),
    @(89, 208),
    @(90, 234),
    @(91, 258),
    @(92, 306),
    @(93, 321),
    @(94, 330),
    @(95, 333),
    @(96, 349),
    @(97, 357),
    @(98, 365),
    @(99, 366),
    @(100, 369),
    @(101, 377),
    @(102, 382),
    @(103, 388),
    @(104, 391),
    @(105, 400),
    @(106, 409),
    @(107, 418),
    @(108, 465),
    @(109, 471),
    @(110, 483),
    @(111, 486),
    @(112, 502),
    @(113, 504),
    @(114, 523),
    @(115, 562),
    @(116, 598),
    @(117, 625),
    @(118, 641),
    @(119, 665),
    @(120, 668),
    @(121, 702),
    @(122, 768),
    @(123, 791),
    @(124, 809),
    @(125, 823),
    @(126, 838),
    @(127, 851),
    @(128, 875),
    @(129, 883),
    @(130, 908),
    @(131, 923),
    @(132, 947),
    @(133, 986),
    @(134, 1052),
    @(135, 1091),
    @(136, 1140),
    @(137, 1233),
    @(138, 1336),
    @(139, 1345),
    @(140, 1398),
    @(141, 1467),
    @(142, 1519),
    @(143, 1573),
    @(144, 1631),
    @(145, 1716),
    @(146, 1858),
    @(147, 1932),
    @(148, 1993),
    @(149, 2088),
    @(150, 2237),
    @(151, 2355),
    @(152, 2386),
    @(153, 2527),
    @(154, 2617),
    @(155, 2731),
    @(156, 2906),
    @(157, 3011),
    @(158, 3082),
    @(159, 3118),
    @(160, 3358),
    @(161, 3408),
    @(162, 3578),
    @(163, 3716),
    @(164, 3829),
    @(165, 4055),
    @(166, 4191),
    @(167, 4294),
    @(168, 4603),
    @(169, 5013),
    @(170, 5189),
    @(171, 5618),
    @(172, 5709),
    @(173, 5838),
    @(174, 6102),
    @(175, 6344),
    @(176, 6786),
    @(177, 7047),
    @(178, 7356),
    @(179, 7785),
    @(180, 8052),
    @(181, 8320),
    @(182, 8794),
    @(183, 8945),
    @(184, 9123),
    @(185, 9481),
    @(186, 9702),
    @(187, 9777)
)

foreach ($pair in $bqUpdates) {
    $row = $pair[0]
    $val = $pair[1]
    $ws.Range("BQ$row").Value = $val
}

# ---------------------------------------------------------------------------
# 2) Scattered corrections elsewhere in the sheet (numeric <-> "NaN" swaps,
#    and a handful of AP/CF off-by-one style recalculations)
# ---------------------------------------------------------------------------
$otherUpdates = @(
    @("I22", "NaN"),
    @("CU26", "NaN"),
    @("DT34", "NaN"),
    @("AP83", 2),
    @("AP88", 7),
    @("CF93", 3),
    @("AP99", 97),
    @("AP100", 105),
    @("AP101", 123),
    @("AP103", "NaN"),
    @("AP104", 251),
    @("AP105", 263),
    @("AP106", 270),
    @("AP107", 276),
    @("AP108", 276),
    @("CF108", "NaN"),
    @("AP109", 295),
    @("AP110", 296),
    @("AP111", 297),
    @("AP112", 298),
    @("AP113", 304),
    @("AP114", 305),
    @("AP115", 308),
    @("AP116", 309),
    @("AP117", 316),
    @("AP118", 316),
    @("AP119", 318),
    @("AP120", 319),
    @("AP121", 323),
    @("AP122", 326),
    @("AP123", 331),
    @("AP124", 334),
    @("AP125", 337),
    @("AP126", 337),
    @("AP127", 338),
    @("AP128", 342),
    @("AP129", 347),
    @("AP130", 352),
    @("AP131", 354),
    @("AP132", 355),
    @("AP133", 358),
    @("AP134", 367),
    @("AP135", 373),
    @("AP136", 377),
    @("AP137", 382),
    @("AP138", 384),
    @("AP139", 430),
    @("AP140", 434),
    @("AP141", 442),
    @("AP142", 453),
    @("AP143", 463),
    @("AP144", 475),
    @("AP145", 482),
    @("AP146", 485),
    @("AP147", 488),
    @("AP148", 494),
    @("AP151", 559),
    @("AP152", 576),
    @("AP153", 584),
    @("AP157", 613),
    @("AP158", 615),
    @("AP169", 776),
    @("AP170", 795),
    @("AP171", 803),
    @("AP172", 837)
)

foreach ($pair in $otherUpdates) {
    $addr = $pair[0]
    $val = $pair[1]
    $ws.Range($addr).Value = $val
}

# ---------------------------------------------------------------------------
# 3) Append new row 188 (2020-09-08 snapshot). Clone row 187's formatting
#    (column styles differ across the row: date style in A, "thousands"
#    style across BS:DQ) before writing the new values into it.
# ---------------------------------------------------------------------------
$ws.Rows("187").Copy()
$ws.Rows("188").PasteSpecial(-4104)
$excel.CutCopyMode = 0

$row188 = @(
    @("A188", 44082),
    @("B188", 679513),
    @("C188", 2719),
    @("D188", 90061),
    @("E188", 65294),
    @("F188", 230120),
    @("G188", 26110),
    @("H188", 4527),
    @("I188", 3546),
    @("J188", 7030),
    @("K188", 6647),
    @("L188", 13488),
    @("M188", 3797),
    @("N188", 21407),
    @("O188", 26442),
    @("P188", 6042),
    @("Q188", 6529),
    @("R188", 13322),
    @("S188", 10476),
    @("T188", 15462),
    @("U188", 13012),
    @("V188", 3274),
    @("W188", 1661),
    @("X188", 7257),
    @("Y188", 22295),
    @("Z188", 12654),
    @("AA188", 8672),
    @("AB188", 50986),
    @("AC188", 1322),
    @("AD188", 386),
    @("AE188", 493),
    @("AF188", 455),
    @("AG188", 327),
    @("AH188", 208),
    @("AI188", 369),
    @("AJ188", 1992),
    @("AK188", 3880),
    @("AL188", 36798),
    @("AM188", 7735),
    @("AN188", 2473),
    @("AO188", 39919),
    @("AP188", 1014),
    @("AQ188", 21023),
    @("AR188", 1480),
    @("AS188", 9017),
    @("AT188", 1571),
    @("AU188", 1581),
    @("AV188", 5644),
    @("AW188", 1712),
    @("AX188", 952),
    @("AY188", 2486),
    @("AZ188", 2649),
    @("BA188", 52176),
    @("BB188", 13072),
    @("BC188", 3737),
    @("BD188", 8343),
    @("BE188", 4842),
    @("BF188", 280),
    @("BG188", 1424),
    @("BH188", 2628),
    @("BI188", 734),
    @("BJ188", 2091),
    @("BK188", 8907),
    @("BL188", 8908),
    @("BM188", 9296),
    @("BN188", 13988),
    @("BO188", 1910),
    @("BP188", 840),
    @("BQ188", 10133),
    @("BR188", 8435),
    @("BS188", 9674),
    @("BT188", 1851),
    @("BU188", 1723),
    @("BV188", 4003),
    @("BW188", 3898),
    @("BX188", 1199),
    @("BY188", 4964),
    @("BZ188", 2774),
    @("CA188", 1526),
    @("CB188", 806),
    @("CC188", 2418),
    @("CD188", 2065),
    @("CE188", 1519),
    @("CF188", 1148),
    @("CG188", 5705),
    @("CH188", 1660),
    @("CI188", 1261),
    @("CJ188", 1452),
    @("CK188", 1829),
    @("CL188", 1721),
    @("CM188", 2068),
    @("CN188", 1307),
    @("CO188", 1114),
    @("CP188", 1133),
    @("CQ188", 668),
    @("CR188", 3129),
    @("CS188", 1188),
    @("CT188", 838),
    @("CU188", 842),
    @("CV188", 1519),
    @("CW188", 1371),
    @("CX188", 693),
    @("CY188", 793),
    @("CZ188", 1058),
    @("DA188", 1329),
    @("DB188", 1160),
    @("DC188", 1272),
    @("DD188", 982),
    @("DE188", 320),
    @("DF188", 344),
    @("DG188", 738),
    @("DH188", 657),
    @("DI188", 433),
    @("DJ188", 534),
    @("DK188", 355),
    @("DL188", 632),
    @("DM188", 721),
    @("DN188", 517),
    @("DO188", 482),
    @("DP188", 372),
    @("DQ188", 517),
    @("DR188", 124256),
    @("DS188", 288487),
    @("DT188", 12737),
    @("DU188", 124112),
    @("DV188", 76892),
    @("DW188", 35356),
    @("DX188", 10631)
)

foreach ($pair in $row188) {
    $addr = $pair[0]
    $val = $pair[1]
    $ws.Range($addr).Value = $val
}

# ---------------------------------------------------------------------------
# 4) Move the active selection to the new bottom-right-most cell
# ---------------------------------------------------------------------------
$ws.Range("DX188").Select()
